# Byreddy_LabExam03Grading.xlsx - grading pass "adloori to davuluri completed"
# Fill in Grading Points (column E) and a grading comment (column F) for the
# "Customer Class" section (rows 3-6), and fill in Grading Points for the
# "Product Class" section (rows 10-14).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Customer Class section (S.No 1-4, rows 3-6) ---
# Q1: Author notation -> 0/1, with a grading comment explaining the deduction
$ws.Range("E3").Value = 0
$ws.Range("F3").Value = "(-1) for missing author notation"

# Q2: Constructor -> full 2/2
$ws.Range("E4").Value = 2

# Q3: Getter method -> full 2/2
$ws.Range("E5").Value = 2

# Q4: toString() method -> full 2/2
$ws.Range("E6").Value = 2

# --- Product Class section (S.No 5-9, rows 10-14) ---
# Q5: Two argument constructor -> full 2/2
$ws.Range("E10").Value = 2

# Q6: Getter methods -> full 2/2
$ws.Range("E11").Value = 2

# Q7: hashcode() method -> full 2/2
$ws.Range("E12").Value = 2

# Q8: equals() method -> full 2/2
$ws.Range("E13").Value = 2

# Q9: toString() method -> full 2/2
$ws.Range("E14").Value = 2

# Leave the cursor/selection where the grader last worked (matches the
# author's saved view position).
$ws.Range("E15").Select()
